# Apply the cryptos-list refresh described by the commit diff.
# Only the cells that actually change are touched, matching the unified diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.466.98"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.225.72"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'111.71"
$ws.Range("E5").Value = "  -3.47%  "
$ws.Range("D6").Value = "'291.90"
$ws.Range("E6").Value = "  +10.16%  "
$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "'0.600"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").Value = "'43.69"
$ws.Range("E10").Value = "  -5.35%  "
$ws.Range("D11").Value = "'0.0912"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").Value = "'54.49"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "'8.63"
$ws.Range("E13").Value = "  -5.20%  "
$ws.Range("D14").Value = "'1.04"
$ws.Range("E14").Value = "  +18.68%  "
$ws.Range("D15").Value = "'0.104"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value = "'14.93"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("D17").Value = "2.561.23"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "2.226.33"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "42.495.92"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("D20").Value = "'7.15"
$ws.Range("E20").Value = "  +6.72%  "
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("D22").Value = "'73.62"
$ws.Range("E22").Value = "  +3.08%  "
$ws.Range("E23").Value = "  +15.86%  "
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").Value = "'234.67"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("D26").Value = "'8.90"
$ws.Range("E26").Value = "  -4.37%  "
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("D28").Value = "'11.42"
$ws.Range("E28").Value = "  -6.18%  "
$ws.Range("D29").Value = "'2.20"
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'174.30"
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'37.32"
$ws.Range("E31").Value = "  -7.72%  "
$ws.Range("D32").Value = "'3.13"
$ws.Range("E32").Value = "  -5.11%  "
$ws.Range("D33").Value = "'21.19"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").Value = "'0.0880"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").Value = "'5.64"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D36").Value = "'5.03"
$ws.Range("E36").Value = "  +9.63%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "'0.126"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.18"
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("D39").Value = "'0.0374"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("D40").Value = "'0.104"
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("E41").Value = "  -3.90%  "
$ws.Range("D42").Value = "'71.40"
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("D43").Value = "'0.230"
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "'12.33"
$ws.Range("E45").Value = "  -7.32%  "
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").Value = "'5.37"
$ws.Range("E47").Value = "  -4.53%  "
$ws.Range("E48").Value = "  +3.76%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'8.43"
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'101.21"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.64"
$ws.Range("E51").Value = "  +5.45%  "
